# Update the "carjacking by neighborhood by month" workbook to add data
# through 2022-07-23 (one additional day of carjacking data for July 2022).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the title text in the header row (B1 uses the
# same shared string, so updating the cell value updates the shared string).
$ws.Name = "Through 2022-07-23"
$ws.Range("B1").Value = "July 2022 (through July 23)"

# Increment existing counts for the new day's incidents.
$ws.Range("AD2").Value = 11   # Austin
$ws.Range("P3").Value = 6     # Englewood
$ws.Range("AY4").Value = 2    # Auburn Gresham
$ws.Range("AK5").Value = 4    # Garfield Park
$ws.Range("B6").Value = 7     # Grand Crossing
$ws.Range("I6").Value = 3     # Grand Crossing
$ws.Range("B8").Value = 5     # North Lawndale
$ws.Range("I8").Value = 8     # North Lawndale
$ws.Range("P17").Value = 2    # Washington Heights
$ws.Range("AK23").Value = 2   # South Chicago
$ws.Range("P26").Value = 2    # Little Village
$ws.Range("B29").Value = 8    # Humboldt Park
$ws.Range("I39").Value = 2    # Wicker Park

# New cells that previously had no recorded incidents.
$ws.Range("I11").Value = 1    # Loop
$ws.Range("W15").Value = 1    # Washington Park
$ws.Range("B36").Value = 1    # West Lawn
$ws.Range("I57").Value = 1    # Woodlawn
$ws.Range("P57").Value = 1    # Woodlawn
$ws.Range("B69").Value = 1    # East Side
$ws.Range("I70").Value = 1    # Edgewater
$ws.Range("B71").Value = 1    # Galewood
$ws.Range("B76").Value = 1    # Hermosa
$ws.Range("I76").Value = 1    # Hermosa
